$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-25 Friday" "2025-07-26 Saturday"
Replace-Text "242×8=1936" "924×2=1848"
Replace-Text "455×7=3185" "184×4=736"
Replace-Text "121×9=1089" "128×7=896"
Replace-Text "303×2=606" "251×9=2259"
Replace-Text "808×2=1616" "443×3=1329"
Replace-Text "567×2=1134" "758×4=3032"
Replace-Text "870×9=7830" "304×2=608"
Replace-Text "346×2=692" "311×6=1866"
Replace-Text "664×9=5976" "879×9=7911"
Replace-Text "102×4=408" "835×2=1670"
Replace-Text "576×4=2304" "319×8=2552"
Replace-Text "975×6=5850" "378×7=2646"
Replace-Text "843×5=4215" "194×8=1552"
Replace-Text "877×3=2631" "156×5=780"
Replace-Text "206×3=618" "937×6=5622"
Replace-Text "846×6=5076" "288×6=1728"
Replace-Text "443×7=3101" "909×8=7272"
Replace-Text "937×2=1874" "574×3=1722"
Replace-Text "104×8=832" "689×4=2756"
Replace-Text "423×9=3807" "655×9=5895"
Replace-Text "709×4=2836" "429×5=2145"
Replace-Text "565×8=4520" "793×9=7137"
Replace-Text "428×3=1284" "394×6=2364"
Replace-Text "860×7=6020" "935×5=4675"
Replace-Text "522×6=3132" "234×5=1170"
